# Horarios actualizados Linea 141 - 331
#
# The scraper re-ran and produced a fresh snapshot (scrape timestamp
# 13:51:32, superseding 13:18:40). This refreshes the three worksheets
# (LP1912, LP1912-215, 6203-6173):
#   - updates the "Ultima actualizacion" / "Total filas" banner cells
#   - rewrites the tail of each sheet's data table (rows whose relative
#     order shifted once the newest scrape's rows were merged in, sorted
#     by Hora_Llegada) and appends the brand-new rows scraped this round
#
# Each entry below is @(row, Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada);
# $null means "leave that column alone".

$sheet1Changes = @(
    @(2, "Última actualización: 13:51:32", $null, $null, $null, $null),
    @(3, "Total filas: 243", $null, $null, $null, $null),
    @(95, "08:00:32", "09:22", "16_SANTA ANA", 82, "LP1912"),
    @(96, "07:24:45", "09:22", "17_ROMERO", 118, "LP1912"),
    @(142, "10:13:53", "11:20", "26_HERNANDEZ", 67, "LP1912"),
    @(143, "11:17:39", "11:20", "225_C ROCA-H SUR", 3, "LP1912"),
    @(171, "11:17:39", "12:20", "14_ABASTO", 63, "LP1912"),
    @(172, "10:52:37", "12:20", "215A_EL PATO", 88, "LP1912"),
    @(202, "12:01:11", "13:14", "215D_EL PATO", 73, "LP1912"),
    @(203, "12:50:41", "13:14", "11_ETCHEVERRY", 24, "LP1912"),
    @(216, "12:01:11", "13:46", "17_ROMERO", 105, "LP1912"),
    @(217, "13:18:40", "13:46", "16_SANTA ANA", 28, "LP1912"),
    @(218, "12:35:30", "13:50", "215A_EL PATO", 75, "LP1912"),
    @(219, "12:50:41", "13:50", "11_ETCHEVERRY", 60, "LP1912"),
    @(221, "13:51:32", "13:51", "11_ETCHEVERRY", 0, "LP1912"),
    @(223, "12:01:11", "13:56", "225_GOMEZ", 115, "LP1912"),
    @(224, "12:01:11", "13:57", "16_P MOR-167 Y 521", 116, "LP1912"),
    @(225, "13:18:40", "14:04", "23_HERNANDEZ", 46, "LP1912"),
    @(226, "12:35:30", "14:04", "17_ROMERO", 89, "LP1912"),
    @(227, "12:50:41", "14:06", "16_SANTA ANA", 76, "LP1912"),
    @(228, "13:51:32", "14:12", "15_ABASTO", 21, "LP1912"),
    @(229, "12:35:30", "14:16", "27_EL RETIRO", 101, "LP1912"),
    @(230, "12:50:41", "14:19", "215C_EL PATO", 89, "LP1912"),
    @(231, "12:35:30", "14:20", "215C_EL PATO", 105, "LP1912"),
    @(232, "12:35:30", "14:21", "26_HERNANDEZ", 106, "LP1912"),
    @(233, "13:51:32", "14:30", "16_SANTA ANA", 39, "LP1912"),
    @(234, "12:50:41", "14:44", "14_ABASTO", 114, "LP1912"),
    @(235, "13:18:40", "14:56", "16_P MOR-SANTA ANA", 98, "LP1912"),
    @(236, "13:18:40", "14:58", "215B_EL PATO", 100, "LP1912"),
    @(237, "13:18:40", "15:00", "81_EL PELIGRO", 102, "LP1912"),
    @(238, "13:18:40", "15:05", "10_OLMOS", 107, "LP1912"),
    @(239, "13:51:32", "15:10", "17_ROMERO", 79, "LP1912"),
    @(240, "13:18:40", "15:13", "11_ETCHEVERRY", 115, "LP1912"),
    @(241, "13:51:32", "15:20", "15_ABASTO", 89, "LP1912"),
    @(242, "13:51:32", "15:26", "26_HERNANDEZ", 95, "LP1912"),
    @(243, "13:51:32", "15:32", "84_COLONIA URQUIZA-ESC 49", 101, "LP1912"),
    @(244, "13:51:32", "15:34", "23_HERNANDEZ", 103, "LP1912"),
    @(245, "13:51:32", "15:37", "10_OLMOS", 106, "LP1912"),
    @(246, "13:51:32", "15:38", "215A_EL PATO", 107, "LP1912"),
    @(247, "13:51:32", "15:46", "14_ABASTO", 115, "LP1912"),
    @(248, "13:51:32", "15:46", "16_P MOR-167 Y 521", 115, "LP1912")
)

$sheet2Changes = @(
    @(2, "Última actualización: 13:51:32", $null, $null, $null, $null),
    @(3, "Total filas: 30", $null, $null, $null, $null),
    @(35, "13:51:32", "15:38", "215A_EL PATO", 107, "LP1912")
)

$sheet3Changes = @(
    @(2, "Última actualización: 13:51:32", $null, $null, $null, $null),
    @(3, "Total filas: 37", $null, $null, $null, $null),
    @(41, "13:51:32", "14:53", "215D_LA PLATA", 62, "L6203"),
    @(42, "13:51:32", "15:34", "215A_LA PLATA", 103, "L6173")
)

function Apply-Changes($ws, $changes) {
    foreach ($row in $changes) {
        $r = $row[0]
        if ($null -ne $row[1]) { $ws.Cells.Item($r, 1).Value = $row[1] }
        if ($null -ne $row[2]) { $ws.Cells.Item($r, 2).Value = $row[2] }
        if ($null -ne $row[3]) { $ws.Cells.Item($r, 3).Value = $row[3] }
        if ($null -ne $row[4]) { $ws.Cells.Item($r, 4).Value = $row[4] }
        if ($null -ne $row[5]) { $ws.Cells.Item($r, 5).Value = $row[5] }
    }
}

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

Apply-Changes $ws1 $sheet1Changes
Apply-Changes $ws2 $sheet2Changes
Apply-Changes $ws3 $sheet3Changes
